$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two occurrence records stored on rows 5 and 6 were re-ordered: every
# field that differed between the two records now holds the other
# record's value (fields that were already identical between the two
# records - locality, county, dates, observer, etc. - stay untouched).

$plainCols = @("A", "B", "E", "F", "G", "H", "Q", "R")
foreach ($col in $plainCols) {
    $addr5 = $col + "5"
    $addr6 = $col + "6"
    $v5 = $ws.Range($addr5).Value2
    $v6 = $ws.Range($addr6).Value2
    $ws.Range($addr5).Value2 = $v6
    $ws.Range($addr6).Value2 = $v5
}

# Column I ("Antal") is stored as text throughout the sheet, even though
# the values look numeric, so force a text write (leading apostrophe)
# rather than letting Value2 auto-convert "1" into the number 1, then
# drop the resulting quote-prefix formatting so no style change lingers.
$i5 = $ws.Range("I5").Value2
$i6 = $ws.Range("I6").Value2

$ws.Range("I5").Formula = "'" + $i6
$ws.Range("I5").Style = "Normal"

$ws.Range("I6").Formula = "'" + $i5
$ws.Range("I6").Style = "Normal"

# AO5 ("Murken granlåga med lite bark kvar.") moves to AO6; AO5 becomes empty
$ao5 = $ws.Range("AO5").Value2
$ws.Range("AO6").Value2 = $ao5
$ws.Range("AO5").ClearContents()
